# Update cryptocurrency price (D) and 1h volume change (E) columns
# to the latest scraped values (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.712.13"
$ws.Cells.Item(2, 5).Value = "  +0.17%  "
$ws.Cells.Item(3, 4).Value = "1.601.18"
$ws.Cells.Item(3, 5).Value = "  +0.20%  "
$ws.Cells.Item(4, 4).Value = "'1.01"
$ws.Cells.Item(4, 5).Value = "  +0.25%  "
$ws.Cells.Item(5, 4).Value = "'211.52"
$ws.Cells.Item(5, 5).Value = "  +0.05%  "
$ws.Cells.Item(6, 4).Value = "'0.514"
$ws.Cells.Item(6, 5).Value = "  -0.18%  "
$ws.Cells.Item(7, 5).Value = "  +0.12%  "
$ws.Cells.Item(8, 5).Value = "  +0.07%  "
$ws.Cells.Item(9, 5).Value = "  -0.06%  "
$ws.Cells.Item(10, 4).Value = "'19.68"
$ws.Cells.Item(10, 5).Value = "  +0.66%  "
$ws.Cells.Item(11, 4).Value = "'0.0846"
$ws.Cells.Item(11, 5).Value = "  +0.82%  "
$ws.Cells.Item(12, 4).Value = "1.825.72"
$ws.Cells.Item(12, 5).Value = "  +0.20%  "
$ws.Cells.Item(13, 4).Value = "1.603.68"
$ws.Cells.Item(13, 5).Value = "  -0.14%  "
$ws.Cells.Item(14, 5).Value = "  +0.49%  "
$ws.Cells.Item(15, 5).Value = "  -0.06%  "
$ws.Cells.Item(16, 4).Value = "'65.02"
$ws.Cells.Item(16, 5).Value = "  -0.26%  "
$ws.Cells.Item(17, 4).Value = "0.0₃0738"
$ws.Cells.Item(17, 5).Value = "  +0.18%  "
$ws.Cells.Item(18, 4).Value = "'210.02"
$ws.Cells.Item(18, 5).Value = "  +0.11%  "
$ws.Cells.Item(19, 5).Value = "  +0.11%  "
$ws.Cells.Item(20, 4).Value = "'7.15"
$ws.Cells.Item(20, 5).Value = "  +1.92%  "
$ws.Cells.Item(21, 5).Value = "  -0.26%  "
$ws.Cells.Item(22, 5).Value = "  -3.15%  "
$ws.Cells.Item(23, 4).Value = "'8.98"
$ws.Cells.Item(23, 5).Value = "  +0.02%  "
$ws.Cells.Item(24, 4).Value = "'143.60"
$ws.Cells.Item(24, 5).Value = "  -0.50%  "
$ws.Cells.Item(25, 5).Value = "  +0.44%  "
$ws.Cells.Item(26, 5).Value = "  -0.57%  "
$ws.Cells.Item(27, 5).Value = "  -0.87%  "
$ws.Cells.Item(28, 4).Value = "'15.33"
$ws.Cells.Item(28, 5).Value = "  +0.31%  "
$ws.Cells.Item(29, 5).Value = "  -1.21%  "
$ws.Cells.Item(30, 4).Value = "'1.15"
$ws.Cells.Item(30, 5).Value = "  -0.06%  "
$ws.Cells.Item(31, 5).Value = "  +0.33%  "
$ws.Cells.Item(32, 4).Value = "'2.96"
$ws.Cells.Item(32, 5).Value = "  -0.19%  "
$ws.Cells.Item(33, 4).Value = "1.288.31"
$ws.Cells.Item(33, 5).Value = "  +0.07%  "
$ws.Cells.Item(34, 5).Value = "  +0.71%  "
$ws.Cells.Item(35, 5).Value = "  +0.17%  "
$ws.Cells.Item(36, 5).Value = "  -3.04%  "
$ws.Cells.Item(37, 5).Value = "  +10.42%  "
$ws.Cells.Item(38, 5).Value = "  -0.11%  "
$ws.Cells.Item(39, 4).Value = "'0.830"
$ws.Cells.Item(39, 5).Value = "  -0.53%  "
$ws.Cells.Item(40, 5).Value = "  -2.09%  "
$ws.Cells.Item(41, 5).Value = "  -0.37%  "
$ws.Cells.Item(42, 5).Value = "  -0.08%  "
$ws.Cells.Item(43, 4).Value = "'62.88"
$ws.Cells.Item(43, 5).Value = "  -1.09%  "
$ws.Cells.Item(44, 4).Value = "1.737.37"
$ws.Cells.Item(45, 4).Value = "'90.49"
$ws.Cells.Item(45, 5).Value = "  -0.34%  "
$ws.Cells.Item(47, 5).Value = "  +0.12%  "
$ws.Cells.Item(48, 5).Value = "  +1.37%  "
$ws.Cells.Item(49, 5).Value = "  +0.11%  "
$ws.Cells.Item(50, 4).Value = "'7.44"
$ws.Cells.Item(50, 5).Value = "  +0.72%  "
$ws.Cells.Item(51, 5).Value = "  +1.01%  "
